$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as text (e.g. "573.03"). Excel would
# otherwise auto-detect plain-looking numeric strings as numbers, so
# pre-format those specific cells as Text before writing the new values.
$textPriceCells = @(
    "D5", "D6", "D9", "D10", "D11", "D12", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D33", "D36", "D37", "D39", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D51"
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.465.80'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '3.439.13'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '573.03'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = '158.77'
$ws.Range("E6").Value = '  -1.97%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.438.10'
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("D9").Value = '0.581'
$ws.Range("E9").Value = '  -5.34%  '
$ws.Range("D10").Value = '7.20'
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("D11").Value = '0.121'
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").Value = '4.033.64'
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").Value = '27.47'
$ws.Range("E15").Value = '  -3.10%  '
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  -9.43%  '
$ws.Range("D17").Value = '64.549.34'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '3.458.73'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").Value = '6.15'
$ws.Range("E19").Value = '  -4.68%  '
$ws.Range("D20").Value = '13.72'
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").Value = '378.51'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").Value = '7.91'
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").Value = '72.20'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Value = '0.529'
$ws.Range("E25").Value = '  -3.98%  '
$ws.Range("D26").Value = '0.0000119'
$ws.Range("D27").Value = '9.92'
$ws.Range("E27").Value = '  -2.00%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '1.43'
$ws.Range("E30").Value = '  -5.62%  '
$ws.Range("D31").Value = '6.06'
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("D33").Value = '23.20'
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("E35").Value = '  -3.27%  '
$ws.Range("D36").Value = '161.06'
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("D37").Value = '1.86'
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("D38").Value = '2.880.85'
$ws.Range("E38").Value = '  -3.42%  '
$ws.Range("D39").Value = '0.0743'
$ws.Range("E39").Value = '  -4.48%  '
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("D41").Value = '0.793'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").Value = '4.51'
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '42.87'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0309'
$ws.Range("E45").Value = '  -3.60%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '25.66'
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  +11.04%  '
$ws.Range("D48").Value = '320.44'
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("D49").Value = '1.07'
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("D51").Value = '0.840'
$ws.Range("E51").Value = '  -3.42%  '
